# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change cell B11 on the "Rules" sheet from the text "R40" to the text "1".
# The new value must remain a *text* string (not get auto-converted to the
# number 1), which is what happens if you simply assign Range.Value = "1"
# to a General-formatted cell. To keep it text (and keep the cell's
# existing style untouched) we build the string "1" via a formula in a
# scratch cell, copy it, and paste-special (values only) into B11 - this
# preserves the shared-string ("t=s") storage Excel uses for literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")
$scratch = $ws.Range("Z1")

# Build a literal text value "1" (formula result is a string, not a number).
$scratch.Formula = "=""1"""
$scratch.Copy()

# xlPasteValues = -4163: paste only the value (as text), keep B11's format/style.
$target.PasteSpecial(-4163)

# Remove the scratch cell entirely so it leaves no trace in the sheet.
$scratch.Clear()
